# Rename variable labels in column A (IAMC_variable) of the first worksheet.
# This corresponds to the commit "change naming of variable for the first dictionnary".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A6").Value = 'Price|Final Energy[index]'
$ws.Range("A12").Value = 'Emissions|CO2eq[per capita]'
$ws.Range("A13").Value = 'Emissions[per capita]'
$ws.Range("A15").Value = 'Sea level|Regional Mean rise'
$ws.Range("A31").Value = 'Households|Number'
$ws.Range("A32").Value = 'Households|Number'
$ws.Range("A68").Value = 'GDP|Country[per capita]'
$ws.Range("A69").Value = 'Revenue|Households|Disposable[per capita]'
$ws.Range("A71").Value = 'Value Added|Real'
$ws.Range("A72").Value = 'Value|GDP Deflator'
$ws.Range("A73").Value = 'Value|GDP Real'
$ws.Range("A74").Value = 'Value|Final Demand real|Taxes on Products'
$ws.Range("A75").Value = 'Value|GDP Real|Taxes on Products'
$ws.Range("A80").Value = 'Final Energy[intensity]'
$ws.Range("A83").Value = 'Final Energy|Non-Energy Use'
$ws.Range("A84").Value = 'Final Energy|Non-Energy Use[intensity]'
$ws.Range("A92").Value = 'Final Energy|Flexibility Options'
$ws.Range("A113").Value = 'Primary Energy Intensity|GDP'
$ws.Range("A114").Value = 'Primary Energy Intensity|GDP[Annual change]'
$ws.Range("A119").Value = 'Final Energy|Net[per capita]'
$ws.Range("A120").Value = 'Final Energy[per capita]'
$ws.Range("A122").Value = 'Primary Energy[per capita]'
$ws.Range("A123").Value = 'Final Energy[per capita]'
$ws.Range("A138").Value = 'Emissions[intensity]'
$ws.Range("A139").Value = 'Emissions|CO2eq[intensity]'
$ws.Range("A140").Value = 'Emissions[intensity]'
$ws.Range("A147").Value = 'EROI|Static|PV'
$ws.Range("A148").Value = 'EROI|Final|PV'
$ws.Range("A149").Value = 'EROI|Static'
$ws.Range("A150").Value = 'EROI|Dynamic'
$ws.Range("A152").Value = 'EROI|Static|Global'
$ws.Range("A153").Value = 'EROI|Static'
$ws.Range("A154").Value = 'ESOI|Static'
$ws.Range("A155").Value = 'ESOI|Static|Storage|Pumped Hydropower Plants'
$ws.Range("A161").Value = 'Resource|Consumption[per capita]'

# Update the view state to match (scrolled to row 113, selection at A114).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 113
$ws.Range("A114").Select()
